$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 keeps its existing cell references, but the shared-string text
# they point to is updated in place (relay part swapped out).
$ws.Range("F14").Value = "653-G5NB-1A4-E-DC5 "
$ws.Range("C14").Value = "G5NB-1A4-E-DC5 "
$ws.Range("B14").Value = "Relais 5A 5VDC"

# New BOM lines for rows 15-20 (entered in the same order the original
# author filled them in, so newly-created shared strings line up).
$ws.Range("F15").Value = "80-T491A106K016 "
$ws.Range("C15").Value = "Tantalkondensatoren - fest für SMD 16V 10uF 1206 10%"
$ws.Range("F16").Value = "755-ESR10EZPJ752 "
$ws.Range("B15").Value = "10 uF Kondensator"
$ws.Range("B16").Value = "7k5 Widerstand"
$ws.Range("C16").Value = "Dickfilmwiderstände - SMD 0805 7K5ohm 5% Anti Surge AEC-Q200"
$ws.Range("B17").Value = "1k2 Widerstand"
$ws.Range("F17").Value = "603-RC0805FR-071K2L"
$ws.Range("C17").Value = "Dickfilmwiderstände - SMD 1.2 kOhms 125 mW 0805 1%"
$ws.Range("F18").Value = "78-V8PAM10S-M3/H "
$ws.Range("C18").Value = "Schottky Dioden & Gleichrichter 8A 100V"
$ws.Range("F19").Value = "80-R82EC4100Z370J "
$ws.Range("C19").Value = "Folienkondensatoren 100V 1uF 5% LS=5mm AEC-Q200 "
$ws.Range("B19").Value = "1uF Kondensator"
$ws.Range("F20").Value = "652-SRP1038C-470M "
$ws.Range("C20").Value = "Festinduktivitäten Ind,11x10x3.8mm,47uH 20%,3.2A,Shd,SMD "
$ws.Range("B20").Value = "47uH Spule"
$ws.Range("B18").Value = "Schottky Diode"

# Supplier column ("Mouser"), reusing the existing shared string.
$ws.Range("E15").Value = "Mouser"
$ws.Range("E16").Value = "Mouser"
$ws.Range("E17").Value = "Mouser"
$ws.Range("E18").Value = "Mouser"
$ws.Range("E19").Value = "Mouser"
# Note: E20 intentionally stays blank, matching the source workbook.

# Row heights grow to fit the newly-wrapped descriptions.
$ws.Rows.Item(15).RowHeight = 33
$ws.Rows.Item(16).RowHeight = 49.5
$ws.Rows.Item(17).RowHeight = 33
$ws.Rows.Item(18).RowHeight = 33
$ws.Rows.Item(19).RowHeight = 33
$ws.Rows.Item(20).RowHeight = 49.5

# Update the active selection to match the author's last position.
$ws.Range("B21").Select()
